$d = $word.ActiveDocument
$failures = @()
$ok = $d.Content.Find.Execute("For Retrofit (Alternative 1) 18000 Monte-Carlo simulations were run.", $true, $false, $false, $false, $false, $true, 1, $false, "For Retrofit (Alternative 1) 13000 Monte-Carlo simulations were run.", 2)
if (-not $ok) { $failures += "For Retrofit (Alternative 1) 18000 Monte-Carlo simulations were run." }
$ok = $d.Content.Find.Execute("For New Bridge (Alternative 2) 9000 Monte-Carlo simulations were run.", $true, $false, $false, $false, $false, $true, 1, $false, "For New Bridge (Alternative 2) 12000 Monte-Carlo simulations were run.", 2)
if (-not $ok) { $failures += "For New Bridge (Alternative 2) 9000 Monte-Carlo simulations were run." }
$ok = $d.Content.Find.Execute("The random number seed for these runs was 70819614.", $true, $false, $false, $false, $false, $true, 1, $false, "The random number seed for these runs was 1622842854.", 2)
if (-not $ok) { $failures += "The random number seed for these runs was 70819614." }
$ok = $d.Content.Find.Execute("(2,510,077; 5,078,943)", $true, $false, $false, $false, $false, $true, 1, $false, "(2,515,918; 5,091,539)", 2)
if (-not $ok) { $failures += "(2,510,077; 5,078,943)" }
$ok = $d.Content.Find.Execute("(6,487,640; 11,024,072)", $true, $false, $false, $false, $false, $true, 1, $false, "(6,462,843; 11,034,492)", 2)
if (-not $ok) { $failures += "(6,487,640; 11,024,072)" }
$ok = $d.Content.Find.Execute("(3,456,450; 4,279,445)", $true, $false, $false, $false, $false, $true, 1, $false, "(3,452,587; 4,279,232)", 2)
if (-not $ok) { $failures += "(3,456,450; 4,279,445)" }
$ok = $d.Content.Find.Execute("(7,700,214; 8,989,256)", $true, $false, $false, $false, $false, $true, 1, $false, "(7,709,895; 8,985,296)", 2)
if (-not $ok) { $failures += "(7,700,214; 8,989,256)" }
$ok = $d.Content.Find.Execute("(-1,366,659; 1,335,737)", $true, $false, $false, $false, $false, $true, 1, $false, "(-1,366,948; 1,367,293)", 2)
if (-not $ok) { $failures += "(-1,366,659; 1,335,737)" }
$ok = $d.Content.Find.Execute("(2,086,003; 6,831,803)", $true, $false, $false, $false, $false, $true, 1, $false, "(2,096,156; 6,841,390)", 2)
if (-not $ok) { $failures += "(2,086,003; 6,831,803)" }
$ok = $d.Content.Find.Execute("(0.91; 4.86)", $true, $false, $false, $false, $false, $true, 1, $false, "(0.92; 4.89)", 2)
if (-not $ok) { $failures += "(0.91; 4.86)" }
$ok = $d.Content.Find.Execute("(4.33; 7.31)", $true, $false, $false, $false, $false, $true, 1, $false, "(4.34; 7.32)", 2)
if (-not $ok) { $failures += "(4.33; 7.31)" }
$ok = $d.Content.Find.Execute("(1.30; 2.73)", $true, $false, $false, $false, $false, $true, 1, $false, "(1.30; 2.75)", 2)
if (-not $ok) { $failures += "(1.30; 2.73)" }
$ok = $d.Content.Find.Execute("(-1,898,759; 2,847,042)", $true, $false, $false, $false, $false, $true, 1, $false, "(-1,888,606; 2,856,629)", 2)
if (-not $ok) { $failures += "(-1,898,759; 2,847,042)" }
$ok = $d.Content.Find.Execute("(-0.24; 0.39)", $true, $false, $false, $false, $false, $true, 1, $false, "(-0.25; 0.39)", 2)
if (-not $ok) { $failures += "(-0.24; 0.39)" }
$ok = $d.Content.Find.Execute("(1.61; 4.95)", $true, $false, $false, $false, $false, $true, 1, $false, "(1.60; 4.95)", 2)
if (-not $ok) { $failures += "(1.61; 4.95)" }
$ok = $d.Content.Find.Execute("Retrofit Indirect Loss Reduction: Gaussian distribution with standard deviation of 600000", $true, $false, $false, $false, $false, $true, 1, $false, "Retrofit Indirect Loss Reduction: Gaussian distribution with standard deviation of 600000.00", 2)
if (-not $ok) { $failures += "Retrofit Indirect Loss Reduction: Gaussian distribution with standard deviation of 600000" }
$ok = $d.Content.Find.Execute("Retrofit Response and Recovery: Gaussian distribution with standard deviation of 180000", $true, $false, $false, $false, $false, $true, 1, $false, "Retrofit Response and Recovery: Gaussian distribution with standard deviation of 180000.00", 2)
if (-not $ok) { $failures += "Retrofit Response and Recovery: Gaussian distribution with standard deviation of 180000" }
$ok = $d.Content.Find.Execute("Retrofit Indirect Cost: Triangular distribution with a min of 475000 and a max of 750000", $true, $false, $false, $false, $false, $true, 1, $false, "Retrofit Indirect Cost: Triangular distribution with a min of 475000.00 and a max of 750000.00", 2)
if (-not $ok) { $failures += "Retrofit Indirect Cost: Triangular distribution with a min of 475000 and a max of 750000" }
$ok = $d.Content.Find.Execute("Retrofit Direct Cost: Triangular distribution with a min of 2850000 and a max of 3840000", $true, $false, $false, $false, $false, $true, 1, $false, "Retrofit Direct Cost: Triangular distribution with a min of 2850000.00 and a max of 3840000.00", 2)
if (-not $ok) { $failures += "Retrofit Direct Cost: Triangular distribution with a min of 2850000 and a max of 3840000" }
$ok = $d.Content.Find.Execute("New Bridge Indirect Loss Reduction: Gaussian distribution with standard deviation of 1050000", $true, $false, $false, $false, $false, $true, 1, $false, "New Bridge Indirect Loss Reduction: Gaussian distribution with standard deviation of 1050000.00", 2)
if (-not $ok) { $failures += "New Bridge Indirect Loss Reduction: Gaussian distribution with standard deviation of 1050000" }
$ok = $d.Content.Find.Execute("New Bridge Response and Recovery: Gaussian distribution with standard deviation of 300000", $true, $false, $false, $false, $false, $true, 1, $false, "New Bridge Response and Recovery: Gaussian distribution with standard deviation of 300000.00", 2)
if (-not $ok) { $failures += "New Bridge Response and Recovery: Gaussian distribution with standard deviation of 300000" }
$ok = $d.Content.Find.Execute("Reduced Commute Time: Triangular distribution with a min of 70000 and a max of 115000", $true, $false, $false, $false, $false, $true, 1, $false, "Reduced Commute Time: Triangular distribution with a min of 70000.00 and a max of 115000.00", 2)
if (-not $ok) { $failures += "Reduced Commute Time: Triangular distribution with a min of 70000 and a max of 115000" }
$ok = $d.Content.Find.Execute("New Bridge OMR: Rectangular distribution with a min of 21375 and a max of 30000", $true, $false, $false, $false, $false, $true, 1, $false, "New Bridge OMR: Rectangular distribution with a min of 21375.00 and a max of 30000.00", 2)
if (-not $ok) { $failures += "New Bridge OMR: Rectangular distribution with a min of 21375 and a max of 30000" }
$ok = $d.Content.Find.Execute("Additional Roadwork Indirect Cost: Triangular distribution with a min of 114000 and a max of 144000", $true, $false, $false, $false, $false, $true, 1, $false, "Additional Roadwork Indirect Cost: Triangular distribution with a min of 114000.00 and a max of 144000.00", 2)
if (-not $ok) { $failures += "Additional Roadwork Indirect Cost: Triangular distribution with a min of 114000 and a max of 144000" }
$ok = $d.Content.Find.Execute("Bridge Construction Indirect Cost: Triangular distribution with a min of 166250 and a max of 224000", $true, $false, $false, $false, $false, $true, 1, $false, "Bridge Construction Indirect Cost: Triangular distribution with a min of 166250.00 and a max of 224000.00", 2)
if (-not $ok) { $failures += "Bridge Construction Indirect Cost: Triangular distribution with a min of 166250 and a max of 224000" }
$ok = $d.Content.Find.Execute("Additional Roadwork Direct Cost: Triangular distribution with a min of 2375000 and a max of 3000000", $true, $false, $false, $false, $false, $true, 1, $false, "Additional Roadwork Direct Cost: Triangular distribution with a min of 2375000.00 and a max of 3000000.00", 2)
if (-not $ok) { $failures += "Additional Roadwork Direct Cost: Triangular distribution with a min of 2375000 and a max of 3000000" }
$ok = $d.Content.Find.Execute("Bridge Construction Direct Cost: Triangular distribution with a min of 4037500 and a max of 5440000", $true, $false, $false, $false, $false, $true, 1, $false, "Bridge Construction Direct Cost: Triangular distribution with a min of 4037500.00 and a max of 5440000.00", 2)
if (-not $ok) { $failures += "Bridge Construction Direct Cost: Triangular distribution with a min of 4037500 and a max of 5440000" }
$ok = $d.Content.Find.Execute("Additional Roadwork OMR: Rectangular distribution with a min of 3500 and a max of 4250", $true, $false, $false, $false, $false, $true, 1, $false, "Additional Roadwork OMR: Rectangular distribution with a min of 3500.00 and a max of 4250.00", 2)
if (-not $ok) { $failures += "Additional Roadwork OMR: Rectangular distribution with a min of 3500 and a max of 4250" }
if ($failures.Count -gt 0) { Write-Output ("FAILED: " + ($failures -join " | ")) } else { Write-Output "All replacements applied successfully." }
